# Update recomputed TPM-derived NATMI metrics (ligand/receptor expression,
# specificity scores and edge weights) for the Col3a1-Ddr1 LR-pair sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 22.628972
$ws.Range("H2").Value = 67.886916
$ws.Range("I2").Value = 0.004372730881336598
$ws.Range("J2").Value = 0.004372730881336598
$ws.Range("M2").Value = 0.3360566666666667
$ws.Range("N2").Value = 1.00817
$ws.Range("O2").Value = 0.01570866217798777
$ws.Range("P2").Value = 0.01570866217798777
$ws.Range("Q2").Value = 7.604616900413334
$ws.Range("R2").Value = 68.44155210372
$ws.Range("S2").Value = 0.00006868975221017135
$ws.Range("T2").Value = 0.00006868975221017135

# Row 3
$ws.Range("G3").Value = 22.628972
$ws.Range("H3").Value = 67.886916
$ws.Range("I3").Value = 0.004372730881336598
$ws.Range("J3").Value = 0.004372730881336598
$ws.Range("O3").Value = 0.109316751024163
$ws.Range("P3").Value = 0.1093167510241629
$ws.Range("Q3").Value = 52.920611756584
$ws.Range("R3").Value = 476.285505809256
$ws.Range("S3").Value = 0.0004780127330507416
$ws.Range("T3").Value = 0.0004780127330507415

# Row 4
$ws.Range("G4").Value = 22.628972
$ws.Range("H4").Value = 67.886916
$ws.Range("I4").Value = 0.004372730881336598
$ws.Range("J4").Value = 0.004372730881336598
$ws.Range("M4").Value = 18.491866
$ws.Range("N4").Value = 55.47559800000001
$ws.Range("O4").Value = 0.864385399390831
$ws.Range("P4").Value = 0.864385399390831
$ws.Range("Q4").Value = 418.4519179417521
$ws.Range("R4").Value = 3766.067261475768
$ws.Range("S4").Value = 0.003779724729292756
$ws.Range("T4").Value = 0.003779724729292756

# Row 5
$ws.Range("G5").Value = 22.628972
$ws.Range("H5").Value = 67.886916
$ws.Range("I5").Value = 0.004372730881336598
$ws.Range("J5").Value = 0.004372730881336598
$ws.Range("M5").Value = 0.2265353333333333
$ws.Range("N5").Value = 0.6796059999999999
$ws.Range("O5").Value = 0.01058918740701822
$ws.Range("P5").Value = 0.01058918740701822
$ws.Range("Q5").Value = 5.126261715010666
$ws.Range("R5").Value = 46.13635543509599
$ws.Range("S5").Value = 0.00004630366678292918
$ws.Range("T5").Value = 0.00004630366678292918

# Row 6
$ws.Range("I6").Value = 0.9598063873258337
$ws.Range("J6").Value = 0.9598063873258338
$ws.Range("M6").Value = 0.3360566666666667
$ws.Range("N6").Value = 1.00817
$ws.Range("O6").Value = 0.01570866217798777
$ws.Range("P6").Value = 0.01570866217798777
$ws.Range("Q6").Value = 1669.199425314657
$ws.Range("R6").Value = 15022.79482783191
$ws.Range("S6").Value = 0.01507727429477641
$ws.Range("T6").Value = 0.01507727429477641

# Row 7
$ws.Range("I7").Value = 0.9598063873258337
$ws.Range("J7").Value = 0.9598063873258338
$ws.Range("O7").Value = 0.109316751024163
$ws.Range("P7").Value = 0.1093167510241629
$ws.Range("S7").Value = 0.1049229158746995
$ws.Range("T7").Value = 0.1049229158746995

# Row 8
$ws.Range("I8").Value = 0.9598063873258337
$ws.Range("J8").Value = 0.9598063873258338
$ws.Range("M8").Value = 18.491866
$ws.Range("N8").Value = 55.47559800000001
$ws.Range("O8").Value = 0.864385399390831
$ws.Range("P8").Value = 0.864385399390831
$ws.Range("Q8").Value = 91849.4264861947
$ws.Range("R8").Value = 826644.8383757523
$ws.Range("S8").Value = 0.8296426274465114
$ws.Range("T8").Value = 0.8296426274465115

# Row 9
$ws.Range("I9").Value = 0.9598063873258337
$ws.Range("J9").Value = 0.9598063873258338
$ws.Range("M9").Value = 0.2265353333333333
$ws.Range("N9").Value = 0.6796059999999999
$ws.Range("O9").Value = 0.01058918740701822
$ws.Range("P9").Value = 0.01058918740701822
$ws.Range("Q9").Value = 1125.205019630015
$ws.Range("R9").Value = 10126.84517667014
$ws.Range("S9").Value = 0.01016356970984637
$ws.Range("T9").Value = 0.01016356970984637

# Row 10
$ws.Range("G10").Value = 182.6322073333333
$ws.Range("H10").Value = 547.896622
$ws.Range("I10").Value = 0.0352911079183418
$ws.Range("J10").Value = 0.0352911079183418
$ws.Range("M10").Value = 0.3360566666666667
$ws.Range("N10").Value = 1.00817
$ws.Range("O10").Value = 0.01570866217798777
$ws.Range("P10").Value = 0.01570866217798777
$ws.Range("Q10").Value = 61.37477082241556
$ws.Range("R10").Value = 552.3729374017399
$ws.Range("S10").Value = 0.0005543760921761406
$ws.Range("T10").Value = 0.0005543760921761406

# Row 11
$ws.Range("G11").Value = 182.6322073333333
$ws.Range("H11").Value = 547.896622
$ws.Range("I11").Value = 0.0352911079183418
$ws.Range("J11").Value = 0.0352911079183418
$ws.Range("O11").Value = 0.109316751024163
$ws.Range("P11").Value = 0.1093167510241629
$ws.Range("Q11").Value = 427.1076979782947
$ws.Range("R11").Value = 3843.969281804652
$ws.Range("S11").Value = 0.003857909257676236
$ws.Range("T11").Value = 0.003857909257676236

# Row 12
$ws.Range("G12").Value = 182.6322073333333
$ws.Range("H12").Value = 547.896622
$ws.Range("I12").Value = 0.0352911079183418
$ws.Range("J12").Value = 0.0352911079183418
$ws.Range("M12").Value = 18.491866
$ws.Range("N12").Value = 55.47559800000001
$ws.Range("O12").Value = 0.864385399390831
$ws.Range("P12").Value = 0.864385399390831
$ws.Range("Q12").Value = 3377.210305292218
$ws.Range("R12").Value = 30394.89274762996
$ws.Range("S12").Value = 0.0305051184129408
$ws.Range("T12").Value = 0.0305051184129408

# Row 13
$ws.Range("G13").Value = 182.6322073333333
$ws.Range("H13").Value = 547.896622
$ws.Range("I13").Value = 0.0352911079183418
$ws.Range("J13").Value = 0.0352911079183418
$ws.Range("M13").Value = 0.2265353333333333
$ws.Range("N13").Value = 0.6796059999999999
$ws.Range("O13").Value = 0.01058918740701822
$ws.Range("P13").Value = 0.01058918740701822
$ws.Range("Q13").Value = 41.37264796565911
$ws.Range("R13").Value = 372.3538316909319
$ws.Range("S13").Value = 0.0003737041555486259
$ws.Range("T13").Value = 0.0003737041555486259

# Row 14
$ws.Range("G14").Value = 2.741590666666667
$ws.Range("H14").Value = 8.224772
$ws.Range("I14").Value = 0.00052977387448787
$ws.Range("J14").Value = 0.00052977387448787
$ws.Range("M14").Value = 0.3360566666666667
$ws.Range("N14").Value = 1.00817
$ws.Range("O14").Value = 0.01570866217798777
$ws.Range("P14").Value = 0.01570866217798777
$ws.Range("Q14").Value = 0.9213298208044445
$ws.Range("R14").Value = 8.291968387239999
$ws.Range("S14").Value = 0.000008322038825053645
$ws.Range("T14").Value = 0.000008322038825053645

# Row 15
$ws.Range("G15").Value = 2.741590666666667
$ws.Range("H15").Value = 8.224772
$ws.Range("I15").Value = 0.00052977387448787
$ws.Range("J15").Value = 0.00052977387448787
$ws.Range("O15").Value = 0.109316751024163
$ws.Range("P15").Value = 0.1093167510241629
$ws.Range("Q15").Value = 6.411544248061333
$ws.Range("R15").Value = 57.703898232552
$ws.Range("S15").Value = 0.00005791315873649664
$ws.Range("T15").Value = 0.00005791315873649663

# Row 16
$ws.Range("G16").Value = 2.741590666666667
$ws.Range("H16").Value = 8.224772
$ws.Range("I16").Value = 0.00052977387448787
$ws.Range("J16").Value = 0.00052977387448787
$ws.Range("M16").Value = 18.491866
$ws.Range("N16").Value = 55.47559800000001
$ws.Range("O16").Value = 0.864385399390831
$ws.Range("P16").Value = 0.864385399390831
$ws.Range("Q16").Value = 50.69712723485068
$ws.Range("R16").Value = 456.274145113656
$ws.Range("S16").Value = 0.0004579288020860255
$ws.Range("T16").Value = 0.0004579288020860255

# Row 17
$ws.Range("G17").Value = 2.741590666666667
$ws.Range("H17").Value = 8.224772
$ws.Range("I17").Value = 0.00052977387448787
$ws.Range("J17").Value = 0.00052977387448787
$ws.Range("M17").Value = 0.2265353333333333
$ws.Range("N17").Value = 0.6796059999999999
$ws.Range("O17").Value = 0.01058918740701822
$ws.Range("P17").Value = 0.01058918740701822
$ws.Range("Q17").Value = 0.6210671555368888
$ws.Range("R17").Value = 5.589604399831999
$ws.Range("S17").Value = 0.000005609874840294203
$ws.Range("T17").Value = 0.000005609874840294203
